$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy cell formatting (number formats / styles) from the existing data rows
# down onto the three new rows before writing values into them.
$ws.Range("A2:CU2").Copy() | Out-Null
$ws.Range("A5:CU7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 5: mariaz09@gmail.com
$ws.Range("A5").Value = 42822.991386597219
$row5 = @(
    "mariaz09@gmail.com", "408 Economic Goals", "704 Middle Class and Professional Groups: Positive", "409 Keynesian Demand Management: Positive", "504 Welfare State Expansion", "503 Equality: Positive",
    "504 Welfare State Expansion", "504 Welfare State Expansion", "504 Welfare State Expansion", "503 Equality: Positive", "503 Equality: Positive", "503 Equality: Positive",
    "503 Equality: Positive", "503 Equality: Positive", "409 Keynesian Demand Management: Positive", "409 Keynesian Demand Management: Positive", "409 Keynesian Demand Management: Positive", "409 Keynesian Demand Management: Positive",
    "503 Equality: Positive", "503 Equality: Positive", "503 Equality: Positive", "503 Equality: Positive", "503 Equality: Positive", "503 Equality: Positive",
    "503 Equality: Positive", "503 Equality: Positive", "503 Equality: Positive", "504 Welfare State Expansion", "408 Economic Goals", "305.2 Political Authority: Personal Competence",
    "503 Equality: Positive", "503 Equality: Positive", "408 Economic Goals", "409 Keynesian Demand Management: Positive", "409 Keynesian Demand Management: Positive", "411 Technology and Infrastructure: Positive",
    "503 Equality: Positive", "409 Keynesian Demand Management: Positive", "411 Technology and Infrastructure: Positive", "402 Incentives: Positive", "402 Incentives: Positive", "411 Technology and Infrastructure: Positive",
    "411 Technology and Infrastructure: Positive", "504 Welfare State Expansion", "411 Technology and Infrastructure: Positive", "410 Economic Growth", "410 Economic Growth", "408 Economic Goals",
    "504 Welfare State Expansion", "504 Welfare State Expansion", "504 Welfare State Expansion", "504 Welfare State Expansion", "504 Welfare State Expansion", "503 Equality: Positive",
    "504 Welfare State Expansion", "504 Welfare State Expansion", "408 Economic Goals", "409 Keynesian Demand Management: Positive", "410 Economic Growth", "504 Welfare State Expansion",
    "503 Equality: Positive", "411 Technology and Infrastructure: Positive", "411 Technology and Infrastructure: Positive", "504 Welfare State Expansion", "504 Welfare State Expansion", "402 Incentives: Positive",
    "402 Incentives: Positive", "403 Market Regulation: Positive", "604 Traditional Morality: Negative", "403 Market Regulation: Positive", "000 No meaningful category applies", "412 Controlled Economy: Positive",
    "408 Economic Goals", "403 Market Regulation: Positive", "404 Economic Planning: Positive", "506 Education Expansion", "506 Education Expansion", "503 Equality: Positive",
    "303 Governmental and Administrative Efficiency: Positive", "506 Education Expansion", "506 Education Expansion", "506 Education Expansion", "506 Education Expansion", "506 Education Expansion",
    "506 Education Expansion", "506 Education Expansion", "506 Education Expansion", "503 Equality: Positive", "503 Equality: Positive", "506 Education Expansion",
    "000 No meaningful category applies", "506 Education Expansion", "301 Decentralisation: Positive", "301 Decentralisation: Positive", "506 Education Expansion", "303 Governmental and Administrative Efficiency: Positive",
    "503 Equality: Positive", "000 No meaningful category applies"
)
for ($i = 0; $i -lt $row5.Count; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $row5[$i]
}

# Row 6: adhikariraju1@gmail.com
$ws.Range("A6").Value = 42823.004877025465
$row6 = @(
    "adhikariraju1@gmail.com", "704 Middle Class and Professional Groups: Positive", "704 Middle Class and Professional Groups: Positive", "402 Incentives: Positive", "504 Welfare State Expansion", "505 Welfare State Limitation",
    "504 Welfare State Expansion", "402 Incentives: Positive", "402 Incentives: Positive", "504 Welfare State Expansion", "408 Economic Goals", "504 Welfare State Expansion",
    "704 Middle Class and Professional Groups: Positive", "505 Welfare State Limitation", "402 Incentives: Positive", "401 Free-Market Economy: Positive", "704 Middle Class and Professional Groups: Positive", "504 Welfare State Expansion",
    "412 Controlled Economy: Positive", "412 Controlled Economy: Positive", "503 Equality: Positive", "000 No meaningful category applies", "702 Labour Groups: Negative", "702 Labour Groups: Negative",
    "505 Welfare State Limitation", "505 Welfare State Limitation", "503 Equality: Positive", "504 Welfare State Expansion", "404 Economic Planning: Positive", "702 Labour Groups: Negative",
    "702 Labour Groups: Negative", "602.1 National Way of Life General: Negative", "602.1 National Way of Life General: Negative", "404 Economic Planning: Positive", "404 Economic Planning: Positive", "408 Economic Goals",
    "705 Minority Groups: Positive", "408 Economic Goals", "506 Education Expansion", "402 Incentives: Positive", "402 Incentives: Positive", "506 Education Expansion",
    "506 Education Expansion", "701 Labour Groups: Positive", "601.1 National Way of Life General: Positive", "504 Welfare State Expansion", "505 Welfare State Limitation", "504 Welfare State Expansion",
    "303 Governmental and Administrative Efficiency: Positive", "301 Decentralisation: Positive", "706 Non-Economic Demographic Groups: Positive", "706 Non-Economic Demographic Groups: Positive", "706 Non-Economic Demographic Groups: Positive", "706 Non-Economic Demographic Groups: Positive",
    "706 Non-Economic Demographic Groups: Positive", "601.1 National Way of Life General: Positive", "601.1 National Way of Life General: Positive", "602.1 National Way of Life General: Negative", "408 Economic Goals", "404 Economic Planning: Positive",
    "706 Non-Economic Demographic Groups: Positive", "411 Technology and Infrastructure: Positive", "411 Technology and Infrastructure: Positive", "303 Governmental and Administrative Efficiency: Positive", "303 Governmental and Administrative Efficiency: Positive", "305.6 Rehabilitation and Compensation",
    "402 Incentives: Positive", "402 Incentives: Positive", "603 Traditional Morality: Positive", "411 Technology and Infrastructure: Positive", "305.6 Rehabilitation and Compensation", "411 Technology and Infrastructure: Positive",
    "000 No meaningful category applies", "503 Equality: Positive", "506 Education Expansion", "506 Education Expansion", "507 Education Limitation", "507 Education Limitation",
    "507 Education Limitation", "503 Equality: Positive", "402 Incentives: Positive", "506 Education Expansion", "704 Middle Class and Professional Groups: Positive", "303 Governmental and Administrative Efficiency: Positive",
    "404 Economic Planning: Positive", "506 Education Expansion", "506 Education Expansion", "406 Protectionism: Positive", "406 Protectionism: Positive", "504 Welfare State Expansion",
    "301 Decentralisation: Positive", "303 Governmental and Administrative Efficiency: Positive", "305.6 Rehabilitation and Compensation", "301 Decentralisation: Positive", "506 Education Expansion", "303 Governmental and Administrative Efficiency: Positive",
    "706 Non-Economic Demographic Groups: Positive", "409 Keynesian Demand Management: Positive"
)
for ($i = 0; $i -lt $row6.Count; $i++) {
    $ws.Cells.Item(6, $i + 2).Value = $row6[$i]
}

# Row 7: m.monteiro@mpp.hertie-school.org
$ws.Range("A7").Value = 42823.046213449074
$row7 = @(
    "m.monteiro@mpp.hertie-school.org", "704 Middle Class and Professional Groups: Positive", "704 Middle Class and Professional Groups: Positive", "704 Middle Class and Professional Groups: Positive", "504 Welfare State Expansion", "415 Marxist Analysis: Positive",
    "504 Welfare State Expansion", "305.2 Political Authority: Personal Competence", "504 Welfare State Expansion", "504 Welfare State Expansion", "705 Minority Groups: Positive", "305.2 Political Authority: Personal Competence",
    "503 Equality: Positive", "503 Equality: Positive", "704 Middle Class and Professional Groups: Positive", "402 Incentives: Positive", "503 Equality: Positive", "503 Equality: Positive",
    "503 Equality: Positive", "503 Equality: Positive", "415 Marxist Analysis: Positive", "503 Equality: Positive", "701 Labour Groups: Positive", "701 Labour Groups: Positive",
    "503 Equality: Positive", "503 Equality: Positive", "305.2 Political Authority: Personal Competence", "504 Welfare State Expansion", "701 Labour Groups: Positive", "305.2 Political Authority: Personal Competence",
    "706 Non-Economic Demographic Groups: Positive", "000 No meaningful category applies", "701 Labour Groups: Positive", "701 Labour Groups: Positive", "701 Labour Groups: Positive", "701 Labour Groups: Positive",
    "607.3 Indigenous rights: Positive", "701 Labour Groups: Positive", "402 Incentives: Positive", "402 Incentives: Positive", "402 Incentives: Positive", "701 Labour Groups: Positive",
    "701 Labour Groups: Positive", "504 Welfare State Expansion", "701 Labour Groups: Positive", "000 No meaningful category applies", "000 No meaningful category applies", "000 No meaningful category applies",
    "701 Labour Groups: Positive", "701 Labour Groups: Positive", "701 Labour Groups: Positive", "706 Non-Economic Demographic Groups: Positive", "706 Non-Economic Demographic Groups: Positive", "706 Non-Economic Demographic Groups: Positive",
    "504 Welfare State Expansion", "504 Welfare State Expansion", "504 Welfare State Expansion", "504 Welfare State Expansion", "000 No meaningful category applies", "504 Welfare State Expansion",
    "504 Welfare State Expansion", "504 Welfare State Expansion", "504 Welfare State Expansion", "504 Welfare State Expansion", "706 Non-Economic Demographic Groups: Positive", "504 Welfare State Expansion",
    "402 Incentives: Positive", "504 Welfare State Expansion", "000 No meaningful category applies", "402 Incentives: Positive", "415 Marxist Analysis: Positive", "503 Equality: Positive",
    "504 Welfare State Expansion", "601.1 National Way of Life General: Positive", "303 Governmental and Administrative Efficiency: Positive", "506 Education Expansion", "506 Education Expansion", "503 Equality: Positive",
    "506 Education Expansion", "504 Welfare State Expansion", "506 Education Expansion", "506 Education Expansion", "504 Welfare State Expansion", "303 Governmental and Administrative Efficiency: Positive",
    "000 No meaningful category applies", "102 Foreign Special Relationships: Negative", "503 Equality: Positive", "503 Equality: Positive", "504 Welfare State Expansion", "504 Welfare State Expansion",
    "302 Centralisation: Positive", "303 Governmental and Administrative Efficiency: Positive", "503 Equality: Positive", "301 Decentralisation: Positive", "301 Decentralisation: Positive", "503 Equality: Positive",
    "504 Welfare State Expansion", "303 Governmental and Administrative Efficiency: Positive"
)
for ($i = 0; $i -lt $row7.Count; $i++) {
    $ws.Cells.Item(7, $i + 2).Value = $row7[$i]
}

# Match the row height used throughout the rest of the sheet.
$ws.Range("5:7").RowHeight = 15.75

# Page margins: switch from the metric (2 cm) top/bottom defaults to the
# US-letter 0.75" defaults.
$ws.PageSetup.TopMargin = 54
$ws.PageSetup.BottomMargin = 54
